$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve the current (old) row 7 data by copying it down into a new row 8,
# then overwrite row 7 with the new weekly reading.

# New row 8 = old row 7 values
$ws.Range("A8").Value = 11
$ws.Range("B8").Value = "Vega Monumental Concepción"
$ws.Range("C8").Value = "Bíobío"
$ws.Range("D8").Value = 44749
$ws.Range("D8").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E8").Value = 8
$ws.Range("F8").Value = 100112035
$ws.Range("G8").Value = "Bruselas (repollito)"
$ws.Range("H8").Value = "Sin especificar"
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 90
$ws.Range("K8").Value = 17000
$ws.Range("L8").Value = 18000
$ws.Range("M8").Value = 17556
$ws.Range("N8").Value = "$/malla 15 kilos"
$ws.Range("O8").Value = "Provincia de Quillota"
$ws.Range("P8").Value = 1170
$ws.Range("Q8").Value = 15
$ws.Range("R8").Value = "Hortaliza"

# Update row 7 with the new weekly reading
$ws.Range("D7").Value = 45119
$ws.Range("J7").Value = 50
$ws.Range("K7").Value = 20000
$ws.Range("L7").Value = 20000
$ws.Range("M7").Value = 20000
$ws.Range("P7").Value = 1333
